# HomeJobs.xlsx fix-up:
#   The newest job posting ("UIL_Electronics (Apprenticeship profile)") had
#   been appended at the very bottom of the sheet (row 28) with a wrong /
#   duplicated job-id (25) instead of being inserted near the top of the
#   countdown list with the next sequential id. This script moves that row
#   into its correct place (row 5) and renumbers the S.No/JobID column so
#   the sequence is consistent again, then updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Move the last data row (28) up to row 5 --------------------------
# 1) Insert a blank row at position 5 -> old rows 5-28 shift down to 6-29.
$ws.Rows("5:5").Insert(-4121)

# 2) The row that used to be 28 is now row 29; copy its populated cells
#    (A:L) into the newly inserted row 5.
$src = $ws.Range("A29:L29")
$dst = $ws.Range("A5:L5")
$src.Copy($dst)

# 3) Row 5 previously (before the insert) had data out to column N; clear
#    the now-stale M5:N5 leftovers since the moved-in row only uses A:L.
$ws.Range("M5:N5").Clear()

# 4) Delete the now-empty row 29 left behind by the move.
$ws.Rows("29:29").Delete()

# --- Correct the S.No / JobID column (column A) ------------------------
# Rows 6-28 already carry the right numbers (they travelled with their
# own row content). Only the four rows above the re-inserted row need a
# manual renumber so the countdown sequence (27..1) stays consistent.
$ws.Range("A2").Value = 27
$ws.Range("A3").Value = 26
$ws.Range("A4").Value = 25
$ws.Range("A5").Value = 24

# --- Update the saved view/selection -----------------------------------
$ws.Range("A2").Select()
